$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H88").Value = 4534.6523
$ws_ALC.Range("J88").Value = 6088
$ws_ALC.Range("L88").Value = 6088
$ws_ALC.Range("N88").Value = -6900

$ws_ALC.Range("H91").Value = 4534.6523
$ws_ALC.Range("J91").Value = 6088
$ws_ALC.Range("L91").Value = 6088
$ws_ALC.Range("N91").Value = -8896

$ws_ALC.Range("H112").Value = 315169.78
$ws_ALC.Range("J112").Value = 359658.4
$ws_ALC.Range("L112").Value = 1078975.2
$ws_ALC.Range("N112").Value = -1081191.2

$ws_ALC.Range("H132").Value = 3749.3572
$ws_ALC.Range("I132").Value = 4091.28
$ws_ALC.Range("K132").Value = 12273.84
$ws_ALC.Range("M132").Value = -9743.84

$ws_ALC.Range("H137").Value = 4557.6875
$ws_ALC.Range("I137").Value = 3039.6667
$ws_ALC.Range("J137").Value = 9111.75
$ws_ALC.Range("K137").Value = 9119.000100000001
$ws_ALC.Range("L137").Value = 27335.25
$ws_ALC.Range("M137").Value = -6569.000100000001
$ws_ALC.Range("N137").Value = -32435.25

$ws_ALC.Range("H138").Value = 387745.72
$ws_ALC.Range("I138").Value = 2196.5557
$ws_ALC.Range("K138").Value = 6589.6671
$ws_ALC.Range("M138").Value = -1449.6671

$ws_ARM.Range("H32").Value = 1429.24
$ws_ARM.Range("I32").Value = 1462.0454
$ws_ARM.Range("K32").Value = 1462.0454
$ws_ARM.Range("M32").Value = -1175.0454

$ws_ARM.Range("H88").Value = 10002.667
$ws_ARM.Range("J88").Value = 10002.667
$ws_ARM.Range("L88").Value = 10002.667
$ws_ARM.Range("N88").Value = -10814.667

$ws_ARM.Range("H91").Value = 10002.667
$ws_ARM.Range("J91").Value = 10002.667
$ws_ARM.Range("L91").Value = 10002.667
$ws_ARM.Range("N91").Value = -12810.667

$ws_ARM.Range("H97").Value = 5458.972
$ws_ARM.Range("I97").Value = 2617.6897
$ws_ARM.Range("J97").Value = 17230
$ws_ARM.Range("K97").Value = 2617.6897
$ws_ARM.Range("L97").Value = 17230
$ws_ARM.Range("M97").Value = -2121.6897
$ws_ARM.Range("N97").Value = -18222

$ws_ARM.Range("H122").Value = 2025.8235
$ws_ARM.Range("I122").Value = 2015.5
$ws_ARM.Range("J122").Value = 2037.4375
$ws_ARM.Range("K122").Value = 6046.5
$ws_ARM.Range("L122").Value = 6112.3125
$ws_ARM.Range("M122").Value = -3596.5
$ws_ARM.Range("N122").Value = -11012.3125

$ws_BSM.Range("H20").Value = 836.2105
$ws_BSM.Range("I20").Value = 935.38464
$ws_BSM.Range("K20").Value = 935.38464
$ws_BSM.Range("M20").Value = -688.38464

$ws_BSM.Range("H82").Value = 38602.5
$ws_BSM.Range("I82").Value = 7364.75
$ws_BSM.Range("J82").Value = 59427.668
$ws_BSM.Range("K82").Value = 7364.75
$ws_BSM.Range("L82").Value = 59427.668
$ws_BSM.Range("M82").Value = -6981.75
$ws_BSM.Range("N82").Value = -60193.668

$ws_BSM.Range("H85").Value = 38602.5
$ws_BSM.Range("I85").Value = 7364.75
$ws_BSM.Range("J85").Value = 59427.668
$ws_BSM.Range("K85").Value = 7364.75
$ws_BSM.Range("L85").Value = 59427.668
$ws_BSM.Range("M85").Value = -6038.75
$ws_BSM.Range("N85").Value = -62079.668

$ws_BSM.Range("H100").Value = 41339.2
$ws_BSM.Range("J100").Value = 41339.2
$ws_BSM.Range("L100").Value = 41339.2
$ws_BSM.Range("N100").Value = -43503.2

$ws_BSM.Range("H132").Value = 87099.89999999999
$ws_BSM.Range("J132").Value = 87099.89999999999
$ws_BSM.Range("L132").Value = 87099.89999999999
$ws_BSM.Range("N132").Value = -97219.89999999999

$ws_CRP.Range("H92").Value = 66182
$ws_CRP.Range("J92").Value = 66182
$ws_CRP.Range("L92").Value = 66182
$ws_CRP.Range("N92").Value = -71174

$ws_CRP.Range("H106").Value = 37992
$ws_CRP.Range("J106").Value = 37992
$ws_CRP.Range("L106").Value = 37992
$ws_CRP.Range("N106").Value = -40516

$ws_CUL.Range("H56").Value = 757391.9
$ws_CUL.Range("I56").Value = 757391.9
$ws_CUL.Range("K56").Value = 757391.9
$ws_CUL.Range("M56").Value = -756861.9

$ws_CUL.Range("H92").Value = 240.5
$ws_CUL.Range("I92").Value = 167
$ws_CUL.Range("K92").Value = 501
$ws_CUL.Range("M92").Value = 747

$ws_CUL.Range("H94").Value = 5799.5
$ws_CUL.Range("J94").Value = 5799.5
$ws_CUL.Range("L94").Value = 17398.5
$ws_CUL.Range("N94").Value = -18750.5

$ws_CUL.Range("H134").Value = 3426.5386
$ws_CUL.Range("I134").Value = 3426.5386
$ws_CUL.Range("K134").Value = 10279.6158
$ws_CUL.Range("M134").Value = -5209.6158

$ws_GSM.Range("H29").Value = 12125
$ws_GSM.Range("J29").Value = 14500
$ws_GSM.Range("L29").Value = 14500
$ws_GSM.Range("N29").Value = -15080

$ws_GSM.Range("H80").Value = 17099.4
$ws_GSM.Range("I80").Value = 3999.5
$ws_GSM.Range("K80").Value = 3999.5
$ws_GSM.Range("M80").Value = -3001.5

$ws_GSM.Range("H83").Value = 17099.4
$ws_GSM.Range("I83").Value = 3999.5
$ws_GSM.Range("K83").Value = 19997.5
$ws_GSM.Range("M83").Value = -15005.5

$ws_GSM.Range("H97").Value = 1375.5217
$ws_GSM.Range("I97").Value = 855.5263
$ws_GSM.Range("J97").Value = 3845.5
$ws_GSM.Range("K97").Value = 855.5263
$ws_GSM.Range("L97").Value = 3845.5
$ws_GSM.Range("M97").Value = -359.5263
$ws_GSM.Range("N97").Value = -4837.5

$ws_GSM.Range("H102").Value = 19311.207
$ws_GSM.Range("I102").Value = 1971.375
$ws_GSM.Range("K102").Value = 1971.375
$ws_GSM.Range("M102").Value = -349.375

$ws_GSM.Range("H122").Value = 2514
$ws_GSM.Range("I122").Value = 2516.6365
$ws_GSM.Range("K122").Value = 7549.9095
$ws_GSM.Range("M122").Value = -5099.9095

$ws_LTW.Range("H16").Value = 8886
$ws_LTW.Range("I16").Value = 12056.75
$ws_LTW.Range("K16").Value = 12056.75
$ws_LTW.Range("M16").Value = -11886.75

$ws_LTW.Range("H22").Value = 3126.0557
$ws_LTW.Range("I22").Value = 2058.6667
$ws_LTW.Range("K22").Value = 2058.6667
$ws_LTW.Range("M22").Value = -1763.6667

$ws_LTW.Range("H27").Value = 3126.0557
$ws_LTW.Range("I27").Value = 2058.6667
$ws_LTW.Range("K27").Value = 2058.6667
$ws_LTW.Range("M27").Value = -1951.6667

$ws_LTW.Range("H104").Value = 10454.75
$ws_LTW.Range("J104").Value = 10454.75
$ws_LTW.Range("L104").Value = 10454.75
$ws_LTW.Range("N104").Value = -17442.75

$ws_LTW.Range("H136").Value = 3827.475
$ws_LTW.Range("I136").Value = 3724.2727
$ws_LTW.Range("K136").Value = 11172.8181
$ws_LTW.Range("M136").Value = -8622.8181

$ws_WVR.Range("H64").Value = 79949.5
$ws_WVR.Range("I64").Value = 0
$ws_WVR.Range("K64").Value = 0
$ws_WVR.Range("M64").Value = $null

$ws_WVR.Range("H67").Value = 79949.5
$ws_WVR.Range("I67").Value = 0
$ws_WVR.Range("K67").Value = 0
$ws_WVR.Range("M67").Value = $null

$ws_WVR.Range("H81").Value = 9080.375
$ws_WVR.Range("I81").Value = 7931.6665
$ws_WVR.Range("J81").Value = 9345.462
$ws_WVR.Range("K81").Value = 15863.333
$ws_WVR.Range("L81").Value = 18690.924
$ws_WVR.Range("M81").Value = -14802.333
$ws_WVR.Range("N81").Value = -20812.924

$ws_WVR.Range("H84").Value = 9080.375
$ws_WVR.Range("I84").Value = 7931.6665
$ws_WVR.Range("J84").Value = 9345.462
$ws_WVR.Range("K84").Value = 79316.66500000001
$ws_WVR.Range("L84").Value = 93454.62
$ws_WVR.Range("M84").Value = -74012.66500000001
$ws_WVR.Range("N84").Value = -104062.62

$ws_WVR.Range("H95").Value = 83483.60000000001
$ws_WVR.Range("J95").Value = 79354.5
$ws_WVR.Range("L95").Value = 79354.5
$ws_WVR.Range("N95").Value = -84846.5

$ws_WVR.Range("H101").Value = 38766.832
$ws_WVR.Range("J101").Value = 38766.832
$ws_WVR.Range("L101").Value = 38766.832
$ws_WVR.Range("N101").Value = -45256.832

$ws_WVR.Range("H136").Value = 2822.1
$ws_WVR.Range("I136").Value = 2211.5715
$ws_WVR.Range("K136").Value = 6634.7145
$ws_WVR.Range("M136").Value = -4084.7145

Write-Output "Applied 184 cell edits across 8 sheets"
